$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the existing data rows (2-18) down
# to (3-19) - this is a new subscriber record added at the top of the list.
$ws.Rows("2:2").Insert()

# Force the new row's cells to Text format so the phone number / DDD / date
# are kept as literal strings (matching every other data row) instead of
# being auto-coerced to a number / date serial.
$ws.Range("A2:C2").NumberFormat = "@"

# Populate the newly inserted row 2 with the new subscriber's data.
$ws.Range("A2").Value = "+5521964219027"
$ws.Range("B2").Value = "21"
$ws.Range("C2").Value = "2024-10-18"

# Re-apply the plain data-row formatting (font/fill/border/alignment) from
# the row directly below so the new row matches the sheet's existing style
# instead of inheriting the header row's bold/red formatting from Insert().
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
